$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$snippetIds = @{
    2 = "oKmf"
    3 = "oKmf"
    4 = "oKmf"
    5 = "oKmf"
    6 = "EWWw"
    7 = "fkRB"
    8 = "jw60"
    9 = "n0GD"
    10 = "n0GD"
    11 = "ScZx"
    12 = "0HuQ"
    13 = "Fph3"
    14 = "5hE2"
    15 = "5PeT"
    16 = "9a3L"
    17 = "wPIk"
    18 = "N2nj"
    19 = "OYCT"
    20 = "Rqfr"
    21 = "j2Eo"
    22 = "Vv9k"
    23 = "EDuz"
    24 = "17Dg"
    25 = "1xIF"
    26 = "bGx0"
    27 = "bGx0"
    28 = "ugeS"
    29 = "rMi4"
    30 = "Iu8v"
    31 = "N5Ua"
    32 = "GvGO"
    33 = "90Aj"
    34 = "90Aj"
    35 = "BLqx"
    36 = "Zyoh"
    37 = "mwR8"
    38 = "DiZp"
    39 = "tWdz"
    40 = "tWdz"
    41 = "vlAh"
    42 = "vlAh"
    43 = "O2sZ"
    44 = "O2sZ"
    45 = "Klxp"
    46 = "Klxp"
    47 = "JLYX"
    48 = "JLYX"
    49 = "43cO"
    50 = "43cO"
    51 = "ZNMU"
    52 = "tclJ"
    53 = "tclJ"
    54 = "DWQB"
    55 = "DWQB"
    56 = "DPh5"
    57 = "DPh5"
    58 = "DPh5"
    59 = "kukZ"
    60 = "kukZ"
    61 = "EDLn"
    62 = "VVAn"
    63 = "4xa0"
    64 = "wVpR"
}

foreach ($row in $snippetIds.Keys) {
    $ws.Cells.Item($row, 8).Value = $snippetIds[$row]
}
